$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p125v_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p125v_1</id>", 2)
$d.Content.Find.Execute("<id>p125v_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p125v_2</id>", 2)
$d.Content.Find.Execute("<id>p125v_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p125v_3</id>", 2)
